$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(
    2,
    1,
    0,
    1,
    1,
    0,
    0,
    3,
    2,
    1,
    1,
    2,
    1,
    2,
    0,
    1,
    2,
    1,
    2,
    2,
    0,
    1,
    1,
    1,
    2,
    0,
    1,
    2,
    1,
    0,
    0,
    1,
    0,
    2,
    1,
    1,
    1,
    1,
    1,
    0,
    1,
    2,
    1,
    1,
    2,
    2,
    1,
    2,
    2,
    1,
    2,
    2,
    1,
    3,
    1,
    0,
    2,
    0,
    1,
    2,
    3,
    3,
    2,
    2,
    1,
    1,
    2,
    2,
    1,
    2,
    2,
    2,
    1,
    1,
    0,
    1,
    0,
    2,
    1,
    1,
    2,
    1,
    2,
    2,
    2,
    2
)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}

Write-Host "Updated K column (G2:G87) with" $kValues.Length "values"
